$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.116.59"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "1.907.76"
$ws.Range("E3").Value = "  +5.41%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'252.49"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5086"
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("D8").Value = "'45.18"
$ws.Range("E8").Value = "  +4.52%  "
$ws.Range("D9").Value = "'0.3032"
$ws.Range("E9").Value = "  +8.95%  "
$ws.Range("D10").Value = "'0.06818"
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("D11").Value = "1.907.77"
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("D12").Value = "'17.32"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "'0.07328"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "'0.6915"
$ws.Range("E14").Value = "  +6.90%  "
$ws.Range("D15").Value = "'86.86"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "'4.916"
$ws.Range("E16").Value = "  +4.70%  "
$ws.Range("D17").Value = "30.114.90"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "'0.000008185"
$ws.Range("E18").Value = "  +11.71%  "
$ws.Range("D19").Value = "'0.9983"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'13.08"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("D21").Value = "2.152.97"
$ws.Range("E21").Value = "  +4.97%  "
$ws.Range("D22").Value = "'0.9982"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'4.821"
$ws.Range("E23").Value = "  +5.08%  "
$ws.Range("E24").Value = "  +7.26%  "
$ws.Range("D25").Value = "'9.284"
$ws.Range("E25").Value = "  +4.58%  "
$ws.Range("D26").Value = "'147.50"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").Value = "'135.31"
$ws.Range("E27").Value = "  +4.62%  "
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").Value = "'2.000"
$ws.Range("E29").Value = "  +5.82%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'4.283"
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").Value = "'0.08848"
$ws.Range("E32").Value = "  +5.94%  "
$ws.Range("D33").Value = "'4.011"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("D34").Value = "'0.05058"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'1.142"
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("D36").Value = "'0.7235"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'2.815"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").Value = "'2.271"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "'0.9649"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'0.01694"
$ws.Range("E41").Value = "  +6.19%  "
$ws.Range("D42").Value = "'6.154"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'0.4309"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("D44").Value = "'105.02"
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'7.635"
$ws.Range("E46").Value = "  +6.57%  "
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("D48").Value = "'0.05746"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").Value = "'33.18"
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("D50").Value = "'8.426"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "'0.3824"
$ws.Range("E51").Value = "  +5.08%  "
